$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "latitude_API" (row 6) and "longitude_API" (row 7) rows entirely.
# Rows below shift up automatically (old row8 -> row6, old row9 -> row7, old row10 -> row8).
$ws.Range("A6:A7").EntireRow.Delete()

# Row 6 is now "depth" - update its definition text to the new data-product wording,
# wrap the text, and grow the row to fit (matches the taller row in the target sheet).
$ws.Range("B6").Value = "Data product depth of sample below sea surface, for underway samples depth of ship's intake. URI http://vocab.nerc.ac.uk/collection/P09/current/DEPH/"
$ws.Range("B6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 31.2

# Update the active selection shown in the sheet view.
$ws.Range("B18:B19").Select() | Out-Null
